$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Global amplitude" (and related) values in column D per the diff
$ws.Cells.Item(11, 4).Value = 5.0
$ws.Cells.Item(22, 4).Value = 4.0
$ws.Cells.Item(33, 4).Value = 4.0
$ws.Cells.Item(44, 4).Value = 4.0
$ws.Cells.Item(55, 4).Value = 4.0
$ws.Cells.Item(66, 4).Value = 4.0
$ws.Cells.Item(77, 4).Value = 4.0
$ws.Cells.Item(88, 4).Value = 5.0
$ws.Cells.Item(99, 4).Value = 3.5
$ws.Cells.Item(110, 4).Value = 3.5
$ws.Cells.Item(121, 4).Value = 3.5
$ws.Cells.Item(132, 4).Value = 3.5
$ws.Cells.Item(143, 4).Value = 3.5
$ws.Cells.Item(154, 4).Value = 4.0
$ws.Cells.Item(165, 4).Value = 5.0
$ws.Cells.Item(176, 4).Value = 4.0
$ws.Cells.Item(187, 4).Value = 3.0
$ws.Cells.Item(198, 4).Value = 4.0
$ws.Cells.Item(209, 4).Value = 3.0
$ws.Cells.Item(220, 4).Value = 2.5
$ws.Cells.Item(231, 4).Value = 5.0
$ws.Cells.Item(242, 4).Value = 5.0
$ws.Cells.Item(253, 4).Value = 4.0
$ws.Cells.Item(264, 4).Value = 4.0
$ws.Cells.Item(275, 4).Value = 5.0
$ws.Cells.Item(286, 4).Value = 2.5
$ws.Cells.Item(297, 4).Value = 2.5
$ws.Cells.Item(308, 4).Value = 5.0
$ws.Cells.Item(466, 4).Value = 35.812
$ws.Cells.Item(467, 4).Value = 10.3335
$ws.Cells.Item(468, 4).Value = 22.0
$ws.Cells.Item(473, 4).Value = 5.0
$ws.Cells.Item(477, 4).Value = 35.432
$ws.Cells.Item(478, 4).Value = 10.6525
$ws.Cells.Item(479, 4).Value = 28.0
$ws.Cells.Item(484, 4).Value = 4.0
$ws.Cells.Item(488, 4).Value = 35.212
$ws.Cells.Item(489, 4).Value = 10.974
$ws.Cells.Item(490, 4).Value = 26.0
$ws.Cells.Item(495, 4).Value = 4.0
$ws.Cells.Item(499, 4).Value = 34.872
$ws.Cells.Item(500, 4).Value = 11.33
$ws.Cells.Item(501, 4).Value = 31.0
$ws.Cells.Item(506, 4).Value = 4.0
$ws.Cells.Item(510, 4).Value = 34.54
$ws.Cells.Item(511, 4).Value = 11.714
$ws.Cells.Item(512, 4).Value = 30.0
$ws.Cells.Item(517, 4).Value = 4.0
$ws.Cells.Item(521, 4).Value = 34.372
$ws.Cells.Item(522, 4).Value = 12.1275
$ws.Cells.Item(523, 4).Value = 28.0
$ws.Cells.Item(528, 4).Value = 4.0
$ws.Cells.Item(532, 4).Value = 34.072
$ws.Cells.Item(533, 4).Value = 12.552
$ws.Cells.Item(534, 4).Value = 36.0
$ws.Cells.Item(539, 4).Value = 5.0

# Restore view/selection state (scroll position + active cell)
$ws.Range("G269").Select()
